# Update "想去人数" (want-to-go count) figures to the freshly scraped values.
# These came from a re-run of the scraper (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 5188
$wsExpo.Range("F5").Value = 5188
$wsExpo.Range("F12").Value = 8659
$wsExpo.Range("F13").Value = 8659
$wsExpo.Range("F32").Value = 7062
$wsExpo.Range("F49").Value = 3272

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5188
$wsAll.Range("F4").Value = 5188
$wsAll.Range("F11").Value = 8659
$wsAll.Range("F12").Value = 8659
$wsAll.Range("F33").Value = 7062
$wsAll.Range("F47").Value = 3273

$wb.Save()
